$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TM")

$ws.Range("D8").Value = 265590800
$ws.Range("E8").Value = 249478600
$ws.Range("F8").Value = 256764200
$ws.Range("G8").Value = 246200100
$ws.Range("H8").Value = 232254900
$ws.Range("I8").Value = 199460300
$ws.Range("J8").Value = 167996200
$ws.Range("D9").Value = 215957900
$ws.Range("E9").Value = 205518400
$ws.Range("F9").Value = 204353400
$ws.Range("G9").Value = 197448800
$ws.Range("H9").Value = 188042300
$ws.Range("I9").Value = 168514600
$ws.Range("J9").Value = 148152600
$ws.Range("D10").Value = 49632800
$ws.Range("E10").Value = 43960200
$ws.Range("F10").Value = 52410800
$ws.Range("G10").Value = 48751300
$ws.Range("H10").Value = 44212600
$ws.Range("I10").Value = 30945700
$ws.Range("J10").Value = 19843600
$ws.Range("D17").Value = 243896000
$ws.Range("E17").Value = 231449500
$ws.Range("F17").Value = 230964300
$ws.Range("G17").Value = 221335000
$ws.Range("H17").Value = 211534200
$ws.Range("I17").Value = 187519500
$ws.Range("J17").Value = 164781400
$ws.Range("D18").Value = 21694800
$ws.Range("E18").Value = 18029100
$ws.Range("F18").Value = 25799900
$ws.Range("G18").Value = 24865100
$ws.Range("H18").Value = 20720700
$ws.Range("I18").Value = 11940800
$ws.Range("J18").Value = 3214900
$ws.Range("D20").Value = 2243300
$ws.Range("E20").Value = 2068400
$ws.Range("F20").Value = 1489900
$ws.Range("G20").Value = 1492800
$ws.Range("H20").Value = 1524100
$ws.Range("I20").Value = 955800
$ws.Range("J20").Value = 905500
$ws.Range("D21").Value = 39631100
$ws.Range("E21").Value = 34676600
$ws.Range("F21").Value = 42003600
$ws.Range("G21").Value = 39110000
$ws.Range("H21").Value = 33565000
$ws.Range("I21").Value = 22897800
$ws.Range("J21").Value = 13784200
$ws.Range("D22").Value = 249400
$ws.Range("E22").Value = 265400
$ws.Range("F22").Value = 320000
$ws.Range("G22").Value = 206800
$ws.Range("H22").Value = 177500
$ws.Range("I22").Value = 207600
$ws.Range("J22").Value = 207200
$ws.Range("D23").Value = 23688700
$ws.Range("E23").Value = 19832200
$ws.Range("F23").Value = 26969800
$ws.Range("G23").Value = 26151200
$ws.Range("H23").Value = 22067400
$ws.Range("I23").Value = 12689000
$ws.Range("J23").Value = 3913200
$ws.Range("D24").Value = 4559800
$ws.Range("E24").Value = 5685300
$ws.Range("F24").Value = 7939600
$ws.Range("G24").Value = 8077000
$ws.Range("H24").Value = 6941000
$ws.Range("I24").Value = 4987200
$ws.Range("J24").Value = 2370900
$ws.Range("D26").Value = 19128800
$ws.Range("E26").Value = 14146900
$ws.Range("F26").Value = 19030200
$ws.Range("G26").Value = 18074200
$ws.Range("H26").Value = 15126400
$ws.Range("I26").Value = 7701700
$ws.Range("J26").Value = 1542200
$ws.Range("D27").Value = 22434500
$ws.Range("E27").Value = 16464700
$ws.Range("F27").Value = 20851700
$ws.Range("G27").Value = 19647000
$ws.Range("H27").Value = 16481000
$ws.Range("I27").Value = 8698000
$ws.Range("J27").Value = 2563400
$ws.Range("D32").Value = -2243300
$ws.Range("E32").Value = -2068400
$ws.Range("F32").Value = -1489900
$ws.Range("G32").Value = -1492800
$ws.Range("H32").Value = -1524100
$ws.Range("I32").Value = -955800
$ws.Range("J32").Value = -905500
$ws.Range("D33").Value = 22434500
$ws.Range("E33").Value = 16464700
$ws.Range("F33").Value = 20851700
$ws.Range("G33").Value = 19647000
$ws.Range("H33").Value = 16481000
$ws.Range("I33").Value = 8698000
$ws.Range("J33").Value = 2563400
$ws.Range("D35").Value = 22434500
$ws.Range("E35").Value = 16464700
$ws.Range("F35").Value = 20851700
$ws.Range("G35").Value = 19647000
$ws.Range("H35").Value = 16481000
$ws.Range("I35").Value = 8698000
$ws.Range("J35").Value = 2563400
$ws.Range("D41").Value = 27592500
$ws.Range("E41").Value = 27075500
$ws.Range("F41").Value = 26572400
$ws.Range("G41").Value = 20652400
$ws.Range("H41").Value = 18452200
$ws.Range("I41").Value = 15533400
$ws.Range("J41").Value = 15180000
$ws.Range("D42").Value = 24133200
$ws.Range("E42").Value = 26254400
$ws.Range("F42").Value = 22992500
$ws.Range("G42").Value = 26500000
$ws.Range("H42").Value = 20132800
$ws.Range("I42").Value = 14033400
$ws.Range("J42").Value = 11402800
$ws.Range("D43").Value = 81877100
$ws.Range("E43").Value = 79095100
$ws.Range("F43").Value = 75612700
$ws.Range("G43").Value = 79545000
$ws.Range("H43").Value = 72467800
$ws.Range("I43").Value = 67999000
$ws.Range("J43").Value = 58970400
$ws.Range("D44").Value = 22959700
$ws.Range("E44").Value = 21593100
$ws.Range("F44").Value = 18636100
$ws.Range("G44").Value = 19324100
$ws.Range("H44").Value = 17128100
$ws.Range("I44").Value = 15510700
$ws.Range("J44").Value = 14665400
$ws.Range("D45").Value = 7537400
$ws.Range("E45").Value = 7198500
$ws.Range("F45").Value = 20800600
$ws.Range("G45").Value = 16123500
$ws.Range("H45").Value = 13907100
$ws.Range("I45").Value = 11538900
$ws.Range("J45").Value = 11165000
$ws.Range("D46").Value = 164100000
$ws.Range("E46").Value = 161217000
$ws.Range("F46").Value = 164614000
$ws.Range("G46").Value = 162145000
$ws.Range("H46").Value = 142088000
$ws.Range("I46").Value = 124615000
$ws.Range("J46").Value = 111384000
$ws.Range("D47").Value = 186824000
$ws.Range("E47").Value = 176849000
$ws.Range("F47").Value = 169476000
$ws.Range("G47").Value = 176925000
$ws.Range("H47").Value = 156772000
$ws.Range("I47").Value = 129067000
$ws.Range("J47").Value = 105167000
$ws.Range("D48").Value = 92819800
$ws.Range("E48").Value = 92181900
$ws.Range("F48").Value = 88053400
$ws.Range("G48").Value = 84033300
$ws.Range("H48").Value = 69077300
$ws.Range("I48").Value = 61935200
$ws.Range("J48").Value = 56367800
$ws.Range("D52").Value = 11042400
$ws.Range("E52").Value = 10453900
$ws.Range("F52").Value = 6601600
$ws.Range("G52").Value = 8374600
$ws.Range("H52").Value = 6656900
$ws.Range("I52").Value = 5151100
$ws.Range("J52").Value = 4166100
$ws.Range("D54").Value = 454787000
$ws.Range("E54").Value = 440702000
$ws.Range("F54").Value = 428746000
$ws.Range("G54").Value = 431478000
$ws.Range("H54").Value = 374595000
$ws.Range("I54").Value = 320769000
$ws.Range("J54").Value = 277085000
$ws.Range("D57").Value = 23383400
$ws.Range("E57").Value = 23200100
$ws.Range("F57").Value = 21601200
$ws.Range("G57").Value = 21791700
$ws.Range("H57").Value = 20007500
$ws.Range("I57").Value = 19108600
$ws.Range("J57").Value = 20273000
$ws.Range("D58").Value = 84444400
$ws.Range("E58").Value = 83566900
$ws.Range("F58").Value = 77030600
$ws.Range("G58").Value = 81030000
$ws.Range("H58").Value = 70335600
$ws.Range("I58").Value = 61417400
$ws.Range("J58").Value = 53908000
$ws.Range("D59").Value = 53056200
$ws.Range("E59").Value = 49796400
$ws.Range("F59").Value = 47133200
$ws.Range("G59").Value = 45719000
$ws.Range("H59").Value = 42370300
$ws.Range("I59").Value = 36203300
$ws.Range("J59").Value = 32324500
$ws.Range("D60").Value = 160884000
$ws.Range("E60").Value = 156563000
$ws.Range("F60").Value = 145765000
$ws.Range("G60").Value = 148541000
$ws.Range("H60").Value = 132713000
$ws.Range("I60").Value = 116729000
$ws.Range("J60").Value = 106505000
$ws.Range("D61").Value = 90457600
$ws.Range("E61").Value = 89600800
$ws.Range("F61").Value = 88339500
$ws.Range("G61").Value = 90530100
$ws.Range("H61").Value = 77264100
$ws.Range("I61").Value = 66333900
$ws.Range("J61").Value = 54622200
$ws.Range("D62").Value = 23349500
$ws.Range("E62").Value = 25770100
$ws.Range("F62").Value = 31123700
$ws.Range("G62").Value = 32875000
$ws.Range("H62").Value = 27037700
$ws.Range("I62").Value = 22239500
$ws.Range("J62").Value = 15916100
$ws.Range("D66").Value = 280966000
$ws.Range("E66").Value = 277976000
$ws.Range("F66").Value = 273016000
$ws.Range("G66").Value = 279713000
$ws.Range("H66").Value = 243794000
$ws.Range("I66").Value = 210951000
$ws.Range("J66").Value = 181710000
$ws.Range("D70").Value = 4447400
$ws.Range("E70").Value = 4392300
$ws.Range("F70").Value = 4337200
$ws.Range("D72").Value = 176040000
$ws.Range("E72").Value = 159114000
$ws.Range("F72").Value = 151820000
$ws.Range("G72").Value = 140951000
$ws.Range("H72").Value = 127611000
$ws.Range("I72").Value = 114710000
$ws.Range("J72").Value = 107730000
$ws.Range("D76").Value = 169373000
$ws.Range("E76").Value = 158334000
$ws.Range("F76").Value = 151392000
$ws.Range("G76").Value = 151765000
$ws.Range("H76").Value = 130801000
$ws.Range("I76").Value = 109818000
$ws.Range("J76").Value = 95374400
$ws.Range("D81").Value = 22434500
$ws.Range("E81").Value = 16464700
$ws.Range("F81").Value = 20851700
$ws.Range("G81").Value = 19647000
$ws.Range("H81").Value = 16481000
$ws.Range("I81").Value = 8698000
$ws.Range("J81").Value = 2563400
$ws.Range("D83").Value = 15675700
$ws.Range("E83").Value = 14563000
$ws.Range("F83").Value = 14697600
$ws.Range("G83").Value = 12738000
$ws.Range("H83").Value = 11307700
$ws.Range("I83").Value = 9990200
$ws.Range("J83").Value = 9653200
$ws.Range("D89").Value = 38058500
$ws.Range("E89").Value = 30864700
$ws.Range("F89").Value = 40326100
$ws.Range("G89").Value = 33319200
$ws.Range("H89").Value = 32960200
$ws.Range("I89").Value = 22159900
$ws.Range("J89").Value = 13130000
$ws.Range("D91").Value = -32532300
$ws.Range("E91").Value = -32014600
$ws.Range("F91").Value = -36695300
$ws.Range("G91").Value = -30352400
$ws.Range("H91").Value = -24215400
$ws.Range("I91").Value = -17846300
$ws.Range("J91").Value = -13850000
$ws.Range("D94").Value = -33087200
$ws.Range("E94").Value = -26848200
$ws.Range("F94").Value = -28770200
$ws.Range("G94").Value = -34473900
$ws.Range("H94").Value = -39199700
$ws.Range("I94").Value = -27366900
$ws.Range("J94").Value = -13041600
$ws.Range("D96").Value = -5667100
$ws.Range("E96").Value = -5769100
$ws.Range("F96").Value = -6381800
$ws.Range("G96").Value = -5016600
$ws.Range("H96").Value = -3580100
$ws.Range("I96").Value = -1717700
$ws.Range("J96").Value = -1417300
$ws.Range("D100").Value = -4060200
$ws.Range("E100").Value = -3391500
$ws.Range("F100").Value = -3829100
$ws.Range("G100").Value = 2766600
$ws.Range("H100").Value = 8312100
$ws.Range("I100").Value = 4314300
$ws.Range("J100").Value = -3212300
$ws.Range("D101").Value = -394000
$ws.Range("E101").Value = -121900
$ws.Range("F101").Value = -1806800
$ws.Range("G101").Value = 588300
$ws.Range("H101").Value = 846200
$ws.Range("I101").Value = 1246200
$ws.Range("J101").Value = -505700
$ws.Range("D102").Value = 517000
$ws.Range("E102").Value = 503000
$ws.Range("F102").Value = 5920000
$ws.Range("G102").Value = 2200200
$ws.Range("H102").Value = 2918800
$ws.Range("I102").Value = 353400
$ws.Range("J102").Value = -3629600
